$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2 and R2 to whole numbers
$ws.Range("Q2").Value = 470097
$ws.Range("R2").Value = 7039164

# Remove the time values in Z2 (Starttid) and AB2 (Sluttid) entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
